# Update "Datos actualizados" timestamp in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 31 de Agosto de 2020 a las 07:10"

# Helper: update a data row (columns A..H) in one shot using a Range write
function Set-Row($r, $a, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("A$r").Value = $a
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
}

# Rows whose country stays the same but whose statistics were refreshed
Set-Row 6   "India"                         3621245 2076 2774801 781827 0 0 64617
Set-Row 19  "Pakistan"                      295849  213  280682  8873   0 6 6294
Set-Row 59  "Kirguistan"                    43898   78   38459   4381   0 0 1058
Set-Row 62  "Uzbekistan"                    41528   104  38819   2394   0 2 315
Set-Row 124 "Tailandia"                     3412    1    3252    102    0 0 58
Set-Row 172 "Islas Turcas y Caicos"         507     2    200     304    0 0 3
Set-Row 175 "San Martin (Parte Holandesa)"  463     3    179     267    0 0 17
Set-Row 196 "Curazao"                       68      1    35      32     0 0 1

# "Jamaica" moves up in the ranking (right after Estonia, row 136) with fresh
# numbers; the countries that used to occupy rows 137-139 (Guinea-Bisau,
# Bahamas, Benin) shift down one row, keeping their own (unchanged) figures.
Set-Row 137 "Jamaica"       2357 244 890  1446 0 1 21
Set-Row 138 "Guinea-Bisau"  2205 0   1127 1044 0 0 34
Set-Row 139 "Bahamas"       2167 0   782  1335 0 0 50
Set-Row 140 "Benin"         2145 0   1738 367  0 0 40

# "Butan" moves up in the ranking (right after Camboya, row 185) with fresh
# numbers; the countries that used to occupy rows 186-187 (San Martin (Parte
# Francesa), Islas Caimanes) shift down one row, keeping their own figures.
Set-Row 186 "Butan"                          224 29 140 84  0 0 0
Set-Row 187 "San Martin (Parte Francesa)"    213 0  79  129 0 0 5
Set-Row 188 "Islas Caimanes"                 205 0  202 2   0 0 1
